$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 18377.8
$ws.Range("I9").Value = 45211
$ws.Range("J9").Value = 489
$ws.Range("K9").Value = 45211
$ws.Range("L9").Value = 489
$ws.Range("M9").Value = -45042
$ws.Range("N9").Value = -827
# Row 11
$ws.Range("H11").Value = 1165681.8
$ws.Range("I11").Value = 1165681.8
$ws.Range("K11").Value = 1165681.8
$ws.Range("M11").Value = -1165541.8
# Row 28
$ws.Range("H28").Value = 4851.25
$ws.Range("I28").Value = 6402
$ws.Range("J28").Value = 2266.6667
$ws.Range("K28").Value = 6402
$ws.Range("L28").Value = 2266.6667
$ws.Range("M28").Value = -5917
$ws.Range("N28").Value = -3236.6667
# Row 32
$ws.Range("H32").Value = 16666.555
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20652
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 40
$ws.Range("H40").Value = 1313.75
$ws.Range("I40").Value = 1230.1428
$ws.Range("J40").Value = 1899
$ws.Range("K40").Value = 1230.1428
$ws.Range("L40").Value = 1899
$ws.Range("M40").Value = -1055.1428
$ws.Range("N40").Value = -2249
# Row 51
$ws.Range("H51").Value = 13891889
$ws.Range("I51").Value = 19233076
$ws.Range("J51").Value = 4799.8
$ws.Range("K51").Value = 19233076
$ws.Range("L51").Value = 4799.8
$ws.Range("M51").Value = -19232592
$ws.Range("N51").Value = -5767.8
# Row 107
$ws.Range("H107").Value = 2027.1212
$ws.Range("I107").Value = 1180.6666
$ws.Range("J107").Value = 4284.3335
$ws.Range("K107").Value = 1180.6666
$ws.Range("L107").Value = 4284.3335
$ws.Range("M107").Value = 739.3334
$ws.Range("N107").Value = -8124.3335
# Row 127
$ws.Range("H127").Value = 3648.875
$ws.Range("I127").Value = 3165.1667
$ws.Range("J127").Value = 5100
$ws.Range("K127").Value = 9495.500100000001
$ws.Range("L127").Value = 15300
$ws.Range("M127").Value = -4535.500100000001
$ws.Range("N127").Value = -25220

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 48924.91
$ws.Range("I32").Value = 34215.312
$ws.Range("K32").Value = 34215.312
$ws.Range("M32").Value = -33928.312
# Row 74
$ws.Range("H74").Value = 55570056
$ws.Range("I74").Value = 4415
$ws.Range("K74").Value = 4415
$ws.Range("M74").Value = -3541
# Row 77
$ws.Range("H77").Value = 55570056
$ws.Range("I77").Value = 4415
$ws.Range("K77").Value = 22075
$ws.Range("M77").Value = -17707

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4885.5557
$ws.Range("I20").Value = 4516.0527
$ws.Range("K20").Value = 4516.0527
$ws.Range("M20").Value = -4269.0527
# Row 134
$ws.Range("H134").Value = 3009.2856
$ws.Range("I134").Value = 2856.818
$ws.Range("K134").Value = 8570.454000000002
$ws.Range("M134").Value = -6035.454000000002

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 4500
$ws.Range("I17").Value = 4500
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -4326
$ws.Range("N17").ClearContents()
# Row 25
$ws.Range("H25").Value = 12537
$ws.Range("I25").Value = 5075
$ws.Range("K25").Value = 5075
$ws.Range("M25").Value = -4901
# Row 31
$ws.Range("H31").Value = 6727784
$ws.Range("I31").Value = 3166644
$ws.Range("J31").Value = 15630633
$ws.Range("K31").Value = 3166644
$ws.Range("L31").Value = 15630633
$ws.Range("M31").Value = -3166349
$ws.Range("N31").Value = -15631223
# Row 34
$ws.Range("H34").Value = 6727784
$ws.Range("I34").Value = 3166644
$ws.Range("J34").Value = 15630633
$ws.Range("K34").Value = 3166644
$ws.Range("L34").Value = 15630633
$ws.Range("M34").Value = -3166442
$ws.Range("N34").Value = -15631037
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 99
$ws.Range("H99").Value = 3013.25
$ws.Range("I99").Value = 2899.3076
$ws.Range("J99").Value = 3224.8572
$ws.Range("K99").Value = 2899.3076
$ws.Range("L99").Value = 3224.8572
$ws.Range("M99").Value = -1401.3076
$ws.Range("N99").Value = -6220.8572
# Row 107
$ws.Range("H107").Value = 840.5
$ws.Range("I107").Value = 971.5714
$ws.Range("J107").Value = 534.6667
$ws.Range("K107").Value = 971.5714
$ws.Range("L107").Value = 534.6667
$ws.Range("M107").Value = 948.4286
$ws.Range("N107").Value = -4374.6667
# Row 126
$ws.Range("H126").Value = 3013.25
$ws.Range("I126").Value = 2899.3076
$ws.Range("J126").Value = 3224.8572
$ws.Range("K126").Value = 8697.9228
$ws.Range("L126").Value = 9674.571599999999
$ws.Range("M126").Value = -6227.9228
$ws.Range("N126").Value = -14614.5716

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 2825
$ws.Range("J39").Value = 3190
$ws.Range("L39").Value = 9570
$ws.Range("N39").Value = -10158

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 37500
$ws.Range("J39").Value = 37500
$ws.Range("L39").Value = 37500
$ws.Range("N39").Value = -38564
# Row 113
$ws.Range("H113").Value = 5237.5
$ws.Range("I113").Value = 5225
$ws.Range("J113").Value = 5250
$ws.Range("K113").Value = 5225
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = -3055
$ws.Range("N113").Value = -9590
# Row 132
$ws.Range("H132").Value = 41922.57
$ws.Range("I132").Value = 41922.57
$ws.Range("K132").Value = 125767.71
$ws.Range("M132").Value = -123237.71

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H23").Value = 3249999.8
$ws.Range("J23").Value = 5500000
$ws.Range("L23").Value = 5500000
$ws.Range("N23").Value = -5500460
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
# Row 46
$ws.Range("H46").Value = 1487.5358
$ws.Range("I46").Value = 848.6087
$ws.Range("K46").Value = 848.6087
$ws.Range("M46").Value = -660.6087
# Row 47
$ws.Range("H47").Value = 18966.3
$ws.Range("J47").Value = 18966.3
$ws.Range("L47").Value = 18966.3
$ws.Range("N47").Value = -19946.3
# Row 52
$ws.Range("H52").Value = 18966.3
$ws.Range("J52").Value = 18966.3
$ws.Range("L52").Value = 18966.3
$ws.Range("N52").Value = -19432.3
# Row 93
$ws.Range("H93").Value = 1113277
$ws.Range("I93").Value = 1113277
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1113277
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1112029
$ws.Range("N93").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 17522.123
$ws.Range("I62").Value = 17518.441
$ws.Range("K62").Value = 17518.441
$ws.Range("M62").Value = -16894.441
# Row 65
$ws.Range("H65").Value = 17522.123
$ws.Range("I65").Value = 17518.441
$ws.Range("K65").Value = 87592.20499999999
$ws.Range("M65").Value = -84472.20499999999
# Row 96
$ws.Range("H96").Value = 2149.75
$ws.Range("I96").Value = 1333
$ws.Range("K96").Value = 1333
$ws.Range("M96").Value = 40
# Row 113
$ws.Range("H113").Value = 1505.4348
$ws.Range("I113").Value = 1707.4706
$ws.Range("J113").Value = 933
$ws.Range("K113").Value = 5122.4118
$ws.Range("L113").Value = 2799
$ws.Range("M113").Value = -2952.4118
$ws.Range("N113").Value = -7139
# Row 122
$ws.Range("H122").Value = 8866.5
$ws.Range("I122").Value = 4049.875
$ws.Range("J122").Value = 18499.75
$ws.Range("K122").Value = 12149.625
$ws.Range("L122").Value = 55499.25
$ws.Range("M122").Value = -9699.625
$ws.Range("N122").Value = -60399.25
# Row 132
$ws.Range("H132").Value = 4572.488
$ws.Range("I132").Value = 4651.8
$ws.Range("K132").Value = 13955.4
$ws.Range("M132").Value = -11425.4
# Row 136
$ws.Range("H136").Value = 1594.2858
$ws.Range("I136").Value = 1476.6957
$ws.Range("J136").Value = 2135.2
$ws.Range("K136").Value = 4430.0871
$ws.Range("L136").Value = 6405.599999999999
$ws.Range("M136").Value = -1880.0871
$ws.Range("N136").Value = -11505.6
